# commit: "add new api thumbs_up"
# The "grimmdb" schema sheet gets one new row describing a new
# "thumbs_up" column (mirrors the existing "interested" row's type/bytes
# /reference/constraint shape, with its own description + comment).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 53

$ws.Range("C$row").Value = "thumbs_up"
$ws.Range("D$row").Value = "TINYINT"
$ws.Range("E$row").Value = 1
$ws.Range("F$row").Value = "NA"
$ws.Range("G$row").Value = "NOT NULL"
$ws.Range("H$row").Value = "是否喜欢"
$ws.Range("I$row").Value = "类似为点赞，默认值为0"

# Leave the new row's cell as the active selection, matching the
# author's saved cursor position after typing the last value.
$ws.Range("I$row").Select()
